$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The author replaced the text values in F4/F6 (which had been pasted in
# as text like "234.33" / " 229.72" with trailing newlines, hence the
# taller row height) with clean numeric values, matching F5's style/format.
$ws.Range("F4").Value = 234.33332999999999
$ws.Range("F6").Value = 229.72

# With the embedded-newline text gone, those rows no longer need the extra
# height that was needed to show two lines of wrapped text - they go back
# to the sheet's normal row height (same as every other row).
$ws.Rows.Item(4).RowHeight = 15.75
$ws.Rows.Item(6).RowHeight = 15.75

# Cursor/selection moved from P8 to F8.
$ws.Range("F8").Select()
